# Applies the "fad.xlsx" update:
#  - renames several column headers (row 1) in Sheet1
#  - updates the GDP-like figures in column C (rows 2-50)
#  - flips the "Colony" flag (column AL) to 1 for rows 6, 16 and 50

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) renames -------------------------------------------------
$ws.Range("C1").Value  = "GDP"
$ws.Range("E1").Value  = "Budget_Previous_Year"
$ws.Range("F1").Value  = "LatinAmerica"
$ws.Range("G1").Value  = "Africa"
$ws.Range("H1").Value  = "Confessional"
$ws.Range("I1").Value  = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) value updates, rows 2-50 --------------------------------
$cValues = @(
    @{ Row = 2; Val = 2898.942214704482 },
    @{ Row = 3; Val = 1904.346464968814 },
    @{ Row = 4; Val = 1503.870423231357 },
    @{ Row = 5; Val = 2983.242707849043 },
    @{ Row = 6; Val = 5555.389721901988 },
    @{ Row = 7; Val = 5082.354756663512 },
    @{ Row = 8; Val = 4633.590358399045 },
    @{ Row = 9; Val = 1955.461557360978 },
    @{ Row = 10; Val = 6336.709213679884 },
    @{ Row = 11; Val = 4355.934938677345 },
    @{ Row = 12; Val = 2965.153206179127 },
    @{ Row = 13; Val = 1939.33862702996 },
    @{ Row = 14; Val = 1577.487171555845 },
    @{ Row = 15; Val = 3083.80337578809 },
    @{ Row = 16; Val = 5660.517066940175 },
    @{ Row = 17; Val = 6947.966251196303 },
    @{ Row = 18; Val = 5360.226632400601 },
    @{ Row = 19; Val = 4921.848409120176 },
    @{ Row = 20; Val = 2024.117324382548 },
    @{ Row = 21; Val = 6711.616186806423 },
    @{ Row = 22; Val = 6911.59200404802 },
    @{ Row = 23; Val = 5642.578115155247 },
    @{ Row = 24; Val = 5122.180090208862 },
    @{ Row = 25; Val = 3156.723844635973 },
    @{ Row = 26; Val = 6869.640636356445 },
    @{ Row = 27; Val = 1657.651524528445 },
    @{ Row = 28; Val = 7200.731056811853 },
    @{ Row = 29; Val = 5919.20956823756 },
    @{ Row = 30; Val = 5295.682695961288 },
    @{ Row = 31; Val = 3212.740625904757 },
    @{ Row = 32; Val = 6796.064220835697 },
    @{ Row = 33; Val = 1716.389195271215 },
    @{ Row = 34; Val = 2286.013198234259 },
    @{ Row = 35; Val = 7449.08671983612 },
    @{ Row = 36; Val = 5412.131646018807 },
    @{ Row = 37; Val = 3252.634165082374 },
    @{ Row = 38; Val = 1775.027517189621 },
    @{ Row = 39; Val = 5996.49696468919 },
    @{ Row = 40; Val = 2361.056581219794 },
    @{ Row = 41; Val = 7580.275568826287 },
    @{ Row = 42; Val = 5330.539154475424 },
    @{ Row = 43; Val = 3314.741082534716 },
    @{ Row = 44; Val = 1836.014008604312 },
    @{ Row = 45; Val = 6114.227214287786 },
    @{ Row = 46; Val = 7633.969039669125 },
    @{ Row = 47; Val = 5176.058803160127 },
    @{ Row = 48; Val = 3382.563653843273 },
    @{ Row = 49; Val = 1895.214690888655 },
    @{ Row = 50; Val = 7026.178156858586 }
)

foreach ($entry in $cValues) {
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Val
}

# --- Column AL ("Colony") flips to 1 for rows 6, 16, 50 ---------------------
$ws.Range("AL6").Value  = 1
$ws.Range("AL16").Value = 1
$ws.Range("AL50").Value = 1
